$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (duplicate of the Dubai (DSC) / Mumbai Indians match row)
$ws.Range("A4").Value = " Dubai (DSC)"
$ws.Range("B4").Value = " October 31 2020"
$ws.Range("C4").Value = "Mumbai won by 9 wickets (with 34 balls remaining)"
$ws.Range("D4").Value = "Delhi Capitals"
$ws.Range("E4").Value = "Mumbai Indians"
$ws.Range("F4").Value = "Harshal Patel "

# G4:K4 hold numeric-looking text ("5","9","1","0","55.55") that must stay
# stored as text (matches the rest of the sheet, which has no sharedStrings
# table and keeps these as inline "str" cells). Force text via NumberFormat,
# then restore the Normal style so no stray formatting is left behind.
$ws.Range("G4:K4").NumberFormat = "@"
$ws.Range("G4").Value = "5"
$ws.Range("H4").Value = "9"
$ws.Range("I4").Value = "1"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "55.55"
$ws.Range("G4:K4").Style = "Normal"

# Row 5 (duplicate of the Sharjah / Rajasthan Royals match row)
$ws.Range("A5").Value = " Sharjah"
$ws.Range("B5").Value = " October 09 2020"
$ws.Range("C5").Value = "Capitals won by 46 runs"
$ws.Range("D5").Value = "Delhi Capitals"
$ws.Range("E5").Value = "Rajasthan Royals"
$ws.Range("F5").Value = "Harshal Patel "

$ws.Range("G5:K5").NumberFormat = "@"
$ws.Range("G5").Value = "16"
$ws.Range("H5").Value = "15"
$ws.Range("I5").Value = "1"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "106.66"
$ws.Range("G5:K5").Style = "Normal"
